# Updated with Field Name with standard name and Minor release
#
# Semantic edit (confirmed by reconstructing the old/new grids from the
# shared-string indices in the xml diff):
#   1. The whole column F ("RequestProcessingType", which was empty in every
#      data row) is deleted from both worksheets. This shifts every column
#      from G onward one step to the left (G->F, H->G, I->H, ... O->N).
#   2. After the shift, five header cells are renamed to their new
#      "standard" names:
#         H1 (old HTTPAction)       -> "Action"
#         I1 (old ExcludeField)     -> "ExcludeFields"
#         J1 (old HttpStatusCode)   -> "StatusCode"
#         M1 (old security)        -> "Security"
#         N1 (old tags)            -> "Tags"
#   3. The workbook is left with the second sheet active/selected (as in
#      the target bookViews/activeTab), and each sheet keeps its own
#      last-used selection.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the now-unused "RequestProcessingType" column; this shifts all
    # the columns after it left by one, exactly like the target file.
    $ws.Columns("F").Delete()

    # Rename headers to their new standard names (columns after the shift).
    $ws.Range("H1").Value = "Action"
    $ws.Range("I1").Value = "ExcludeFields"
    $ws.Range("J1").Value = "StatusCode"
    $ws.Range("M1").Value = "Security"
    $ws.Range("N1").Value = "Tags"
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Restore each sheet's own selection.
[void]$ws1.Range("J2").Select()
[void]$ws2.Range("O1").Select()

# Second sheet ("API-Testing-Sheet2-Duplicate") is the active tab.
[void]$ws2.Activate()
